# Update scripts wuth new tpm
# Applies updated TPM-derived values to rows 2-6 and removes the now-obsolete
# row 7 (Resolving-Mac target) from the LR-pair sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Target cluster: ECs) ---
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.035285
$ws.Range("H2").Value = 0.07056999999999999
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4562695
$ws.Range("N2").Value = 0.912539
$ws.Range("O2").Value = 0.01609359429837405
$ws.Range("P2").Value = 0.01172153108534722
$ws.Range("Q2").Value = 0.0160994693075
$ws.Range("R2").Value = 0.06439787723
$ws.Range("S2").Value = 0.01609359429837405
$ws.Range("T2").Value = 0.01172153108534722

# --- Row 3 (Target cluster: FAPs) ---
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.035285
$ws.Range("H3").Value = 0.07056999999999999
$ws.Range("O3").Value = 0.7445397209174328
$ws.Range("P3").Value = 0.8134117203287967
$ws.Range("Q3").Value = 0.7448115171099999
$ws.Range("R3").Value = 4.468869102659999
$ws.Range("S3").Value = 0.7445397209174328
$ws.Range("T3").Value = 0.8134117203287967

# --- Row 4 (Target cluster: Inflammatory-Mac) ---
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.035285
$ws.Range("H4").Value = 0.07056999999999999
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.006762666666666667
$ws.Range("N4").Value = 0.020288
$ws.Range("O4").Value = 0.0002385336160064851
$ws.Range("P4").Value = 0.0002605986403425218
$ws.Range("Q4").Value = 0.0002386206933333333
$ws.Range("R4").Value = 0.00143172416
$ws.Range("S4").Value = 0.0002385336160064851
$ws.Range("T4").Value = 0.0002605986403425218

# --- Row 5 (Target cluster: MuSCs) ---
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.035285
$ws.Range("H5").Value = 0.07056999999999999
$ws.Range("M5").Value = 6.745213
$ws.Range("N5").Value = 13.490426
$ws.Range("O5").Value = 0.2379179881147404
$ws.Range("P5").Value = 0.1732840434365834
$ws.Range("Q5").Value = 0.238004840705
$ws.Range("R5").Value = 0.9520193628199999
$ws.Range("S5").Value = 0.2379179881147404
$ws.Range("T5").Value = 0.1732840434365834

# --- Row 6 (Target cluster: Neutrophils) ---
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.035285
$ws.Range("H6").Value = 0.07056999999999999
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.03430933333333333
$ws.Range("N6").Value = 0.102928
$ws.Range("O6").Value = 0.00121016305344615
$ws.Range("P6").Value = 0.00132210650893016
$ws.Range("Q6").Value = 0.001210604826666666
$ws.Range("R6").Value = 0.007263628959999999
$ws.Range("S6").Value = 0.00121016305344615
$ws.Range("T6").Value = 0.00132210650893016

# --- Row 7 (Target cluster: Resolving-Mac) is no longer present in the data ---
$ws.Rows("7:7").Delete()
